$wb = $excel.ActiveWorkbook

# --- Sheet "loai": replace the "dep" category row with "Xe đạp" and give it
#     wrap-text / vertically-centered formatting ---
$wsLoai = $wb.Worksheets.Item("loai")
$rngA4 = $wsLoai.Range("A4")
$rngA4.Value = "Xe đạp"
$rngA4.WrapText = $true
$rngA4.VerticalAlignment = -4108  ## xlCenter

# --- Sheet "chi tiet sp": update product id in row 3 from 2 to 25 ---
$wsChiTiet = $wb.Worksheets.Item("chi tiet sp")
$wsChiTiet.Range("A3").Value = 25

# --- Selections / active sheet bookkeeping ---
$wsChiTiet.Range("H17").Select()
$wsLoai.Range("F9").Select()
